# IP - EX 54 - Graficos avancados (Anexo).xlsx
# "Exercicio" - add a SUB-TOTAL picker (Media/Soma/Contar/...) driven by a
# VLOOKUP into the SUBTOTAL function-number table (S18:T27) plus a
# data-validation dropdown, and wire every P16:P32 + D35:O36 cell to
# SUBTOTAL()/summary formulas instead of the broken array formula that used
# to live in P16.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)   # "EXERCÍCIOS"

# ---------------------------------------------------------------------
# 1) Rebuild the little "SUB-TOTAL" box in M12:P13.
#    Before: M13:N13 ("SUB-TOTAL" label) + O13:P13 (empty) merged, row 12
#    empty. After: the whole box moves up to span rows 12-13, with O12
#    holding the chosen statistic name (dropdown) and P12 holding the
#    VLOOKUP that turns it into a SUBTOTAL function number.
# ---------------------------------------------------------------------

$ws.Range("M13:N13").UnMerge()
$ws.Range("O13:P13").UnMerge()

# Label ("SUB-TOTAL") now lives in M12 (merged M12:N13)
$ws.Range("M13").ClearContents()
$ws.Range("M12").Value = "SUB-TOTAL"
$ws.Range("M12:N13").Merge()
$ws.Range("M12:N13").Font.Bold = $true
$ws.Range("M12:N13").HorizontalAlignment = -4108   # xlCenter
$ws.Range("M12:N13").VerticalAlignment = -4108     # xlCenter
$ws.Range("M12:N13").Interior.Color = $ws.Range("B16").Interior.Color
$ws.Range("M12").Borders.Item(7).LineStyle = 1     # xlEdgeLeft
$ws.Range("M12").Borders.Item(8).LineStyle = 1     # xlEdgeTop
$ws.Range("M12").Borders.Item(10).LineStyle = 1    # xlEdgeRight
$ws.Range("N13").Borders.Item(7).LineStyle = 1     # xlEdgeLeft
$ws.Range("N13").Borders.Item(9).LineStyle = 1     # xlEdgeBottom
$ws.Range("N13").Borders.Item(10).LineStyle = 1    # xlEdgeRight

# Dropdown cell (chosen statistic) in O12 (merged O12:O13)
$ws.Range("O12:O13").Merge()
$ws.Range("O12").Value = "Média"
$ws.Range("O12:O13").Font.Bold = $true
$ws.Range("O12:O13").WrapText = $true
$ws.Range("O12:O13").HorizontalAlignment = -4108
$ws.Range("O12:O13").VerticalAlignment = -4108
$ws.Range("O12:O13").Interior.ColorIndex = -4142
$ws.Range("O12:O13").Interior.Pattern = 1
$ws.Range("O12:O13").Interior.Color = $ws.Range("C2").Interior.Color
$ws.Range("O12").Borders.Item(7).LineStyle = 1
$ws.Range("O12").Borders.Item(8).LineStyle = 1
$ws.Range("O12").Borders.Item(10).LineStyle = 1
$ws.Range("O13").Borders.Item(7).LineStyle = 1
$ws.Range("O13").Borders.Item(9).LineStyle = 1
$ws.Range("O13").Borders.Item(10).LineStyle = 1

# VLOOKUP result (SUBTOTAL function number) in P12 (merged P12:P13)
$ws.Range("P12:P13").Merge()
$ws.Range("P12").Formula = "=VLOOKUP(`$O`$12,`$S`$18:`$T`$27,2,FALSE)"
$ws.Range("P12:P13").Font.Color = $ws.Range("C2").Interior.Color
$ws.Range("P12:P13").HorizontalAlignment = -4108
$ws.Range("P12:P13").VerticalAlignment = -4108
$ws.Range("P12:P13").Interior.Color = $ws.Range("C2").Interior.Color
$ws.Range("P12").Borders.Item(7).LineStyle = 1

# Dropdown list validation on O12 (list of stat names S19:S27)
$ws.Range("O12").Validation.Delete()
$ws.Range("O12").Validation.Add(3, 1, 1, "=`$S`$19:`$S`$27")
$ws.Range("O12").Validation.IgnoreBlank = $true
$ws.Range("O12").Validation.InCellDropdown = $true
$ws.Range("O12").Validation.ShowInput = $true
$ws.Range("O12").Validation.ShowError = $true

# ---------------------------------------------------------------------
# 2) Replace the broken array formula in P16 and fill in P17:P32 with
#    live SUBTOTAL() formulas driven by $P$12.
# ---------------------------------------------------------------------

$ws.Range("P16:P27").Formula = "=SUBTOTAL(`$P`$12,`$D16:`$O16)"
$ws.Range("P28").Formula = "=SUBTOTAL(`$P`$12,`$D28:`$O28)"
$ws.Range("P29:P32").Formula = "=SUBTOTAL(`$P`$12,`$D29:`$O29)"

# ---------------------------------------------------------------------
# 3) Fill in the "% Retrabalho Mensal" (row 35) and "% Retrabalho YTD"
#    (row 36) formulas across D:O.
# ---------------------------------------------------------------------

$ws.Range("D35:O35").Formula = "=SUM(D29:D31)/D32"
$ws.Range("D36").Formula = "=D35"
$ws.Range("E36:O36").Formula = "=D36+E35"

# ---------------------------------------------------------------------
# 4) Column widths for the new M:P box (AutoFit to content) and the
#    selection left by the user on Q42.
# ---------------------------------------------------------------------

$ws.Columns.Item(14).AutoFit()
$ws.Columns.Item(15).AutoFit()

$wb.Application.Calculate()

$ws.Range("Q42").Select()
